$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @(2,2,0.9329097690169021),
    @(2,3,0.05228061303763809),
    @(2,4,0.1156676683935771),
    @(2,5,0.06005454773048058),
    @(2,6,2.075524740652568),
    @(2,8,0.07973214163530429),
    @(2,9,1.68842009516662),
    @(2,11,0.7479508682361029),
    @(2,12,0.2249561431626361),
    @(2,13,0.2387367315177329),
    @(2,14,3.034829201778066),
    @(3,2,0.8966713425653268),
    @(3,3,0.04550233385545255),
    @(3,4,0.1157624700085336),
    @(3,5,0.06020361776108896),
    @(3,6,2.064219328938819),
    @(3,8,0.07973214163530429),
    @(3,9,1.688000342013574),
    @(3,11,0.7065628561795734),
    @(3,12,0.2222859337751899),
    @(3,13,0.2320865127833471),
    @(3,14,3.050060697382932),
    @(4,2,0.874947644033341),
    @(4,3,0.04133518365776467),
    @(4,4,0.1158200975665817),
    @(4,5,0.06030932923153909),
    @(4,6,2.058230358761719),
    @(4,8,0.07973214163530429),
    @(4,9,1.688397463660699),
    @(4,11,0.6815773114576871),
    @(4,12,0.2207555003845627),
    @(4,13,0.22813448093018),
    @(4,14,3.060176053919939),
    @(5,2,0.8662277516443737),
    @(5,3,0.03963563088288424),
    @(5,4,0.1158434391272163),
    @(5,5,0.06035598249805219),
    @(5,6,2.056029287262106),
    @(5,8,0.07973214163530429),
    @(5,9,1.68872397326848),
    @(5,11,0.6715028810065462),
    @(5,12,0.2201592984918292),
    @(5,13,0.2265570411517075),
    @(5,14,3.064489997482482),
    @(6,2,0.8647878407414851),
    @(6,3,0.03935333415734021),
    @(6,4,0.1158473065081118),
    @(6,5,0.06036394541818435),
    @(6,6,2.055678263633098),
    @(6,8,0.07973214163530429),
    @(6,9,1.688788136566615),
    @(6,11,0.6698365165326834),
    @(6,12,0.2200619592099997),
    @(6,13,0.226297105726907),
    @(6,14,3.065217910657928),
    @(7,2,0.874829506970201),
    @(7,3,0.04131226869280624),
    @(7,4,0.1158204129297662),
    @(7,5,0.06030994392691191),
    @(7,6,2.058199704766224),
    @(7,8,0.07973214163530429),
    @(7,9,1.68840120029315),
    @(7,11,0.6814410093192009),
    @(7,12,0.2207473485614742),
    @(7,13,0.2281130731641881),
    @(7,14,3.060233456538818),
    @(8,2,0.9203055264001421),
    @(8,3,0.04994450224798186),
    @(8,4,0.115700478945536),
    @(8,5,0.06010300888559517),
    @(8,6,2.071428879898477),
    @(8,8,0.07973214163530429),
    @(8,9,1.688139452461634),
    @(8,11,0.7335916972909899),
    @(8,12,0.2240128316294587),
    @(8,13,0.23641652452228),
    @(8,14,3.039922598852954),
    @(9,2,1.013660706942687),
    @(9,3,0.06683506623973301),
    @(9,4,0.1154604974818447),
    @(9,5,0.05980936469561016),
    @(9,6,2.104936121109773),
    @(9,8,0.07973214163530429),
    @(9,9,1.692823406491279),
    @(9,11,0.8392504305352304),
    @(9,12,0.231281207882617),
    @(9,13,0.2537399399285576),
    @(9,14,3.006150687704647),
    @(10,2,1.084799269778784),
    @(10,3,0.07922924020277833),
    @(10,4,0.115281007993433),
    @(10,5,0.05966152085894016),
    @(10,6,2.134180169307839),
    @(10,8,0.07973214163530429),
    @(10,9,1.699437477044711),
    @(10,11,0.9189607810573648),
    @(10,12,0.2371483412928512),
    @(10,13,0.2671022574942512),
    @(10,14,2.985034759134209),
    @(11,2,1.117717361209827),
    @(11,3,0.08486599471868317),
    @(11,4,0.1151986118553889),
    @(11,5,0.05960890761050486),
    @(11,6,2.148492139862356),
    @(11,8,0.07973214163530429),
    @(11,9,1.70313648791651),
    @(11,11,0.955679273605341),
    @(11,12,0.2399319355212555),
    @(11,13,0.2733192051427977),
    @(11,14,2.976232028045231),
    @(12,2,1.130262590177949),
    @(12,3,0.08700039087548816),
    @(12,4,0.1151672997287871),
    @(12,5,0.05959108145639469),
    @(12,6,2.154056945733188),
    @(12,8,0.07973214163530429),
    @(12,9,1.704636515718391),
    @(12,11,0.969649570208901),
    @(12,12,0.2410024752727509),
    @(12,13,0.2756932812351991),
    @(12,14,2.973014223189196),
    @(13,2,1.127557202004255),
    @(13,3,0.08654071450061451),
    @(13,4,0.1151740483200321),
    @(13,5,0.0595948274864071),
    @(13,6,2.152852006616115),
    @(13,8,0.07973214163530429),
    @(13,9,1.704309041201441),
    @(13,11,0.9666378909716116),
    @(13,12,0.2407711840584597),
    @(13,13,0.2751810991906964),
    @(13,14,2.973702091713491),
    @(14,2,1.118747865356681),
    @(14,3,0.08504159471047501),
    @(14,4,0.115196038021022),
    @(14,5,0.05960739905446033),
    @(14,6,2.148947049589907),
    @(14,8,0.07973214163530429),
    @(14,9,1.703257905911201),
    @(14,11,0.9568272997890404),
    @(14,12,0.2400196798817831),
    @(14,13,0.2735141241756764),
    @(14,14,2.975964979802967),
    @(15,2,1.113362286660674),
    @(15,3,0.08412332734533834),
    @(15,4,0.1152094928617871),
    @(15,5,0.05961537239891967),
    @(15,6,2.14657406018776),
    @(15,8,0.07973214163530429),
    @(15,9,1.702626987181105),
    @(15,11,0.9508265974652943),
    @(15,12,0.2395615037454775),
    @(15,13,0.2724956380074914),
    @(15,14,2.977366122231459),
    @(16,2,1.082659181757407),
    @(16,3,0.07886084364106694),
    @(16,4,0.1152863774875073),
    @(16,5,0.05966525319139926),
    @(16,6,2.13326515652048),
    @(16,8,0.07973214163530429),
    @(16,9,1.699209630462676),
    @(16,11,0.916570351064621),
    @(16,12,0.2369687306964465),
    @(16,13,0.2666987458258561),
    @(16,14,2.985626207658015),
    @(17,2,1.06396628481798),
    @(17,3,0.07563217656348797),
    @(17,4,0.1153333503884717),
    @(17,5,0.05969959773081168),
    @(17,6,2.12535901523799),
    @(17,8,0.07973214163530429),
    @(17,9,1.697290004949025),
    @(17,11,0.8956724927278685),
    @(17,12,0.2354074822761589),
    @(17,13,0.263177947090611),
    @(17,14,2.990899261574455),
    @(18,2,1.053267045880403),
    @(18,3,0.07377499425182066),
    @(18,4,0.1153602979770465),
    @(18,5,0.05972073042495296),
    @(18,6,2.120906548548362),
    @(18,8,0.07973214163530429),
    @(18,9,1.696250853091229),
    @(18,11,0.8836956949030821),
    @(18,12,0.2345202823212986),
    @(18,13,0.2611659026584263),
    @(18,14,2.994007758045754),
    @(19,2,1.049653477072326),
    @(19,3,0.07314615792314783),
    @(19,4,0.1153694100515708),
    @(19,5,0.05972812261999838),
    @(19,6,2.119415321625326),
    @(19,8,0.07973214163530429),
    @(19,9,1.695910171001046),
    @(19,11,0.8796479602122247),
    @(19,12,0.2342217454879574),
    @(19,13,0.260486897524622),
    @(19,14,2.995073217901989),
    @(20,2,1.065950751730981),
    @(20,3,0.07597588740812),
    @(20,4,0.1153283573103394),
    @(20,5,0.05969579907263167),
    @(20,6,2.12619081100128),
    @(20,8,0.07973214163530429),
    @(20,9,1.697487628358324),
    @(20,11,0.8978926452757889),
    @(20,12,0.2355725634293719),
    @(20,13,0.2635513943133958),
    @(20,14,2.990330113303898),
    @(21,2,1.121333215257607),
    @(21,3,0.08548192513481467),
    @(21,4,0.1151895821371998),
    @(21,5,0.05960364962331255),
    @(21,6,2.150090089173645),
    @(21,8,0.07973214163530429),
    @(21,9,1.703563954741142),
    @(21,11,0.95970712272279),
    @(21,12,0.2402399685592655),
    @(21,13,0.2740032166853652),
    @(21,14,2.97529717638426),
    @(22,2,1.157994261249939),
    @(22,3,0.09169400712809761),
    @(22,4,0.1150982390493596),
    @(22,5,0.05955564581917905),
    @(22,6,2.166555888893512),
    @(22,8,0.07973214163530429),
    @(22,9,1.708113957990719),
    @(22,11,1.000489975969458),
    @(22,12,0.2433862761353538),
    @(22,13,0.2809497852380005),
    @(22,14,2.96614614020146),
    @(23,2,1.138385053603258),
    @(23,3,0.08837853512977745),
    @(23,4,0.1151470506888028),
    @(23,5,0.05958015077753309),
    @(23,6,2.157690309387476),
    @(23,8,0.07973214163530429),
    @(23,9,1.705632564911397),
    @(23,11,0.9786883197213854),
    @(23,12,0.241698267635897),
    @(23,13,0.2772316986358518),
    @(23,14,2.970968519723371),
    @(24,2,1.065053426090145),
    @(24,3,0.07582049879215447),
    @(24,4,0.1153306148590598),
    @(24,5,0.0596975121235932),
    @(24,6,2.125814466821254),
    @(24,8,0.07973214163530429),
    @(24,9,1.697398081998145),
    @(24,11,0.89688879700347),
    @(24,12,0.2354978978984263),
    @(24,13,0.2633825211007021),
    @(24,14,2.990587185536029),
    @(25,2,0.9879580745888461),
    @(25,3,0.06226894508672842),
    @(25,4,0.1155259586163737),
    @(25,5,0.05987684571733176),
    @(25,6,2.09506016370328),
    @(25,8,0.07973214163530429),
    @(25,9,1.690999427525348),
    @(25,11,0.8103022901492238),
    @(25,12,0.2292223546112098),
    @(25,13,0.2489420619134535),
    @(25,14,3.014638037130226)
)

foreach ($item in $changes) {
    $row = $item[0]
    $col = $item[1]
    $val = $item[2]
    $ws.Cells.Item($row, $col).Value = $val
}

Write-Output "Applied $($changes.Count) cell updates"
